# Updates the cryptos price/volume table to the latest scrape (GitHub Actions run).
# Column D ("Price") and E ("Volume(1h)") are plain text cells in the source sheet
# (not numbers), so any value that looks numeric is written back with the cell
# NumberFormat forced to "@" (Text) first -- otherwise Excel's normal type-inference
# on .Value would silently turn e.g. "0.997" or "1.00" into a real number and lose
# the original text formatting (trailing zeros, thousands-style dots, etc.).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "46.415.51"
$ws.Range("E2").Value = "  +1.57%  "
# Row 3
$ws.Range("D3").Value = "2.539.49"
$ws.Range("E3").Value = "  +7.43%  "
# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.15%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "306.28"
$ws.Range("E5").Value = "  +1.92%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.99"
$ws.Range("E6").Value = "  +5.16%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.606"
$ws.Range("E7").Value = "  +6.12%  "
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").Value = "  -0.13%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.570"
$ws.Range("E9").Value = "  +10.81%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.04"
$ws.Range("E10").Value = "  +12.17%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0828"
$ws.Range("E11").Value = "  +3.39%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.95"
$ws.Range("E12").Value = "  +11.54%  "
# Row 13
$ws.Range("D13").Value = "2.925.75"
$ws.Range("E13").Value = "  +7.33%  "
# Row 14
$ws.Range("E14").Value = "  +2.74%  "
# Row 15
$ws.Range("D15").Value = "2.547.95"
$ws.Range("E15").Value = "  +7.29%  "
# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.896"
$ws.Range("E16").Value = "  +10.19%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.94"
# Row 18
$ws.Range("D18").Value = "46.375.12"
$ws.Range("E18").Value = "  +1.42%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.84"
$ws.Range("E19").Value = "  +10.21%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0000100"
$ws.Range("E20").Value = "  +2.40%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.62"
$ws.Range("E21").Value = "  +10.90%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "70.40"
$ws.Range("E22").Value = "  +6.51%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "254.97"
$ws.Range("E23").Value = "  +4.80%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.98"
$ws.Range("E24").Value = "  +5.61%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  +12.44%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  -0.13%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "41.97"
$ws.Range("E27").Value = "  +4.16%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "24.09"
$ws.Range("E28").Value = "  +15.04%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.43"
$ws.Range("E29").Value = "  +7.67%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.26"
$ws.Range("E30").Value = "  +1.22%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.85"
$ws.Range("E31").Value = "  +3.77%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.01"
$ws.Range("E32").Value = "  +11.22%  "
# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("E33").Value = "  +6.44%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0846"
$ws.Range("E34").Value = "  +9.35%  "
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.18"
$ws.Range("E35").Value = "  +18.69%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "150.02"
$ws.Range("E36").Value = "  +3.20%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.117"
$ws.Range("E37").Value = "  +5.94%  "
# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.121"
$ws.Range("E38").Value = "  +4.56%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.57"
$ws.Range("E39").Value = "  +8.33%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.23"
$ws.Range("E40").Value = "  +8.88%  "
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0326"
$ws.Range("E41").Value = "  +9.47%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.56"
$ws.Range("E42").Value = "  +11.98%  "
# Row 43
$ws.Range("D43").Value = "2.015.87"
$ws.Range("E43").Value = "  +8.38%  "
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.998"
$ws.Range("E44").Value = "  +0.05%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "95.02"
$ws.Range("E45").Value = "  +4.86%  "
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "17.41"
$ws.Range("E46").Value = "  +33.95%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.87"
$ws.Range("E47").Value = "  +6.38%  "
# Row 48
$ws.Range("B48").Value = "Algorand"
$ws.Range("C48").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.203"
$ws.Range("E48").Value = "  +10.36%  "
# Row 49
$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.02"
$ws.Range("E49").Value = "  +12.55%  "
# Row 50
$ws.Range("B50").Value = "Aave"
$ws.Range("C50").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "107.18"
$ws.Range("E50").Value = "  +10.80%  "
# Row 51
$ws.Range("D51").Value = "2.791.77"
$ws.Range("E51").Value = "  +7.47%  "
